$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- sheet view: select B2 (drops the B1 topLeftCell pin, moves the selection) ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
$ws.Range("B2").Select() | Out-Null

# --- new Neo4j/Cypher query text added to A2 (cell already carries the wrap-text style) ---
$ws.Range("A2").Value = 'MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN [''Hurthle cell neoplasm (thyroid)''] RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(t.clinical_trial_designation ,'''')as `Trial Code` , coalesce(a.arm_id,'''') As `Arm` , coalesce(a.arm_drug,'''') As `Arm Treatment` , coalesce(c.disease,'''') As Diagnosis , coalesce(c.gender,'''') As Gender , coalesce(c.race,'''') As Race , coalesce(c.ethnicity,'''') As Ethnicity'

# --- row 2 grows tall enough to show the wrapped query text ---
$ws.Rows.Item(2).RowHeight = 87
